$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.650.23'
$ws.Range('E2').Value = '  +0.54%  '

$ws.Range('D3').Value = '1.959.41'
$ws.Range('E3').Value = '  +1.03%  '

$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.75'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.78%  '

$ws.Range('E6').Value = '  +1.12%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '61.62'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +8.16%  '

$ws.Range('E8').Value = '  +0.03%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.376'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.69%  '

$ws.Range('E10').Value = '  -6.54%  '

$ws.Range('E11').Value = '  +0.33%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.27'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +6.25%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.03'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.38%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.834'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.37%  '

$ws.Range('D15').Value = '2.231.95'
$ws.Range('E15').Value = '  +0.33%  '

$ws.Range('E16').Value = '  +3.08%  '

$ws.Range('D17').Value = '1.965.79'
$ws.Range('E17').Value = '  +1.28%  '

$ws.Range('D18').Value = '36.562.61'
$ws.Range('E18').Value = '  +0.45%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.85'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.05%  '

$ws.Range('D20').Value = '0.0₃0853'
$ws.Range('E20').Value = '  -1.24%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '230.30'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.33%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.07'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.92%  '

$ws.Range('E23').Value = '  +0.10%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.45'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +5.29%  '

$ws.Range('E25').Value = '  +3.14%  '

$ws.Range('E26').Value = '  +5.10%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.19'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.39%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '160.36'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.40%  '

$ws.Range('E29').Value = '  +1.08%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.28'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +16.35%  '

$ws.Range('E31').Value = '  +1.50%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.77'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.99%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0616'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.61%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.46'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +7.40%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.55'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +16.11%  '

$ws.Range('E36').Value = '  +0.08%  '

$ws.Range('E37').Value = '  +4.18%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.77'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.15%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.54'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -8.91%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0984'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.56%  '

$ws.Range('E42').Value = '  +2.59%  '

$ws.Range('E43').Value = '  +0.92%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.06'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.09%  '

$ws.Range('D45').Value = '1.369.46'
$ws.Range('E45').Value = '  +2.16%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '88.66'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.91%  '

$ws.Range('E47').Value = '  +1.34%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.14'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.55%  '

$ws.Range('E49').Value = '  +0.47%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '45.75'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.91%  '

$ws.Range('D51').Value = '2.127.27'
$ws.Range('E51').Value = '  +0.51%  '
